$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet shuffle: "demo" -> "rgth" (existing sheet3), then a brand-new
# "pgth" sheet appended at the end (becomes sheet4).
# ---------------------------------------------------------------------
$wsRgth = $wb.Worksheets.Item(3)
$wsRgth.Name = "rgth"

$wsPgth = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsPgth.Name = "pgth"

# ---------------------------------------------------------------------
# dpc sheet (sheet2) content tweaks
# ---------------------------------------------------------------------
$wsDpc = $wb.Worksheets.Item("dpc")
$wsDpc.Range("B2").Value = "COI_Auto"
$wsDpc.Range("A3").Value = "competitor"
$wsDpc.Range("B3").Value = "ak@gmail.com"
$wsDpc.Hyperlinks.Add($wsDpc.Range("B3"), "mailto:ak@gmail.com") | Out-Null
$wsDpc.Range("A4").Value = "company competes"
$wsDpc.Range("A5").Value = "relationship"
$wsDpc.Range("B5").Value = "relationship"
$wsDpc.Range("A6").Value = "additional"
$wsDpc.Range("B6").Value = "Additional information "

# ---------------------------------------------------------------------
# rgth sheet (sheet3) - full rebuild
# ---------------------------------------------------------------------
$wsRgth.Range("A1").Value = "TestCases"
$wsRgth.Range("B1").Value = "Data"
$wsRgth.Range("A2").Value = "Request_Name"
$wsRgth.Range("B2").Value = "Receive Gift Auto"
$wsRgth.Range("A3").Value = "Full_Name"
$wsRgth.Range("B3").Value = "Akshay"
$wsRgth.Range("A4").Value = "Employer"
$wsRgth.Range("B4").Value = "Shreya"
$wsRgth.Range("A5").Value = "Offical_Position"
$wsRgth.Range("B5").Value = "QA"
$wsRgth.Range("A6").Value = "email_address"
$wsRgth.Range("B6").Value = "ak@gmail.com"
$wsRgth.Hyperlinks.Add($wsRgth.Range("B6"), "mailto:ak@gmail.com") | Out-Null
$wsRgth.Range("A7").Value = "Description"
$wsRgth.Range("B7").Value = "Details about the gift, travel or hospitality"
$wsRgth.Range("A8").Value = "Business_Purpose"
$wsRgth.Range("B8").Value = "Business Purpose or Rationale."
$wsRgth.Range("A9").Value = "Monetary_Value"
$wsRgth.Range("B9").Value = 500

$wsRgth.Columns.Item(1).ColumnWidth = 30.140625
$wsRgth.Columns.Item(2).ColumnWidth = 50.85546875

# ---------------------------------------------------------------------
# pgth sheet (sheet4) - full build
# ---------------------------------------------------------------------
$wsPgth.Range("A1").Value = "TestCases"
$wsPgth.Range("B1").Value = "Data"
$wsPgth.Range("A2").Value = "RequestName"
$wsPgth.Range("B2").Value = "Provide Gift Auto"
$wsPgth.Range("A3").Value = "Describe the gift"
$wsPgth.Range("B3").Value = "hospitality in detail, including the business purpose or rationale"
$wsPgth.Range("A4").Value = "Additional Information"
$wsPgth.Range("B4").Value = "upload any additional files to support your approval request (optional)."
$wsPgth.Range("A5").Value = "Sort Description"
$wsPgth.Range("B5").Value = "Category Entertainment"
$wsPgth.Range("A6").Value = "Amount"
$wsPgth.Range("B6").Value = 500
$wsPgth.Range("B6").Style = $wsRgth.Range("B6").Style
$wsPgth.Range("A7").Value = "FullName"
$wsPgth.Range("B7").Value = "akshay"
$wsPgth.Range("A8").Value = "Employer"
$wsPgth.Range("B8").Value = "kapil"
$wsPgth.Range("A9").Value = "Title"
$wsPgth.Range("B9").Value = "Tester"
$wsPgth.Range("A10").Value = "Email"
$wsPgth.Range("B10").Value = "aagg@gmail.com"
$wsPgth.Hyperlinks.Add($wsPgth.Range("B10"), "mailto:aagg@gmail.com") | Out-Null

$wsPgth.Columns.Item(1).ColumnWidth = 24.5703125
$wsPgth.Columns.Item(2).ColumnWidth = 28.85546875

# ---------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------
$wsDpc.Range("D12").Select() | Out-Null
$wsRgth.Range("A15").Select() | Out-Null
$wsPgth.Range("D12").Select() | Out-Null
$wsPgth.Activate()
